# Fruta / hortaliza, semanal
# New weekly Camote (Vega Central Mapocho de Santiago) price data is
# inserted as two new observations right after the existing row 100.
# That pushes the previously-recorded rows 101-166 down to rows 103-168
# (the last two of which become brand-new worksheet rows), and rows 101
# and 102 are populated with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the existing data block (rows 101-166, cols A-R)
# before anything gets overwritten.
$src = $ws.Range("A101:R166").Value2

# --- Step 2: shift that whole block down by two rows, onto rows 103-168.
$ws.Range("A103:R168").Value2 = $src

# New cells created by the shift at the bottom (rows 167-168) need the
# same date number format column D carries everywhere else in the sheet.
$ws.Range("D167:D168").NumberFormat = $ws.Range("D166").NumberFormat

# --- Step 3: write the two brand-new observations into rows 101-102,
# reusing the constant columns (A,B,C,E,F,G,H,O,Q,R) already in place.
$ws.Range("D101").Value2 = 45096
$ws.Range("I101").Value2 = "Primera"
$ws.Range("J101").Value2 = 700
$ws.Range("K101").Value2 = 17000
$ws.Range("L101").Value2 = 19000
$ws.Range("M101").Value2 = 18000
$ws.Range("N101").Value2 = "$/caja 18 kilos"
$ws.Range("P101").Value2 = 1000
$ws.Range("Q101").Value2 = 18

$ws.Range("D102").Value2 = 45096
$ws.Range("I102").Value2 = "Primera"
$ws.Range("J102").Value2 = 520
$ws.Range("K102").Value2 = 16000
$ws.Range("L102").Value2 = 17000
$ws.Range("M102").Value2 = 16500
$ws.Range("N102").Value2 = "$/malla 18 kilos"
$ws.Range("P102").Value2 = 917
$ws.Range("Q102").Value2 = 18
